$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: replace header text cells (A1="version", B1="data") with numeric
# sequence 0,1,2 across A1:C1. C1 is a brand-new cell, so copy A1's
# formatting (bold/border/center style) onto it before writing the value.
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2

# Row 2: the three separate "dataframe dump" rows (old A2/B2, A3/B3, A4/B4)
# collapse into a single row, one text blob per column.
$textA2 = @'
                      Time Driver                LapTime  NumberOfLaps  NumberOfPitStops  ... SpeedI1 SpeedI2 SpeedFL SpeedST IsPersonalBest
0   0 days 01:03:58.671000      1                    NaT             1                 0  ...   274.0   293.0   273.0   239.0          False
1   0 days 01:06:10.776000      1 0 days 00:02:12.105000             2                 0  ...   167.0   107.0   253.0    93.0           True
2   0 days 01:08:21.433000      1 0 days 00:02:10.657000             3                 0  ...   204.0   162.0   289.0   133.0           True
3   0 days 01:09:44.824000      1 0 days 00:01:23.391000             4                 0  ...   280.0   296.0   288.0   303.0           True
4   0 days 01:11:07.928000      1 0 days 00:01:23.104000             5                 0  ...   273.0   294.0   290.0   286.0           True
..                     ...    ...                    ...           ...               ...  ...     ...     ...     ...     ...            ...
990 0 days 02:38:34.078000     11 0 days 00:01:40.500000            54                 3  ...     NaN   163.0   258.0     NaN          False
991 0 days 02:55:25.519000     11                    NaT            55                 3  ...   237.0   150.0     NaN   163.0          False
992 0 days 02:59:01.376000     11                    NaT            56                 4  ...    62.0   142.0   100.0   204.0          False
993 0 days 03:32:15.650000     11                    NaT            57                 4  ...   110.0   199.0     NaN   252.0          False
994 0 days 03:34:58.928000     11                    NaT            58                 5  ...   118.0   137.0   290.0   168.0          False

[995 rows x 18 columns]
'@
$textB2 = @'
                        Time Driver  Position GapToLeader IntervalToPositionAhead
0     0 days 00:00:06.078000      1         1       LAP 1                   LAP 1
1     0 days 01:02:30.972000      1         2       LAP 1                   LAP 1
2     0 days 01:02:31.519000      1         2      +0.624                  +0.624
3     0 days 01:02:35.238000      1         2      +0.672                  +0.672
4     0 days 01:02:38.785000      1         2      +0.823                  +0.823
...                      ...    ...       ...         ...                     ...
21019 0 days 03:34:41.836000     11         6      +3.517                  +0.464
21020 0 days 03:34:49.867000     11         6      +3.728                  +0.543
21021 0 days 03:34:52.336000     11         6      +3.442                  +0.440
21022 0 days 03:34:58.836000     11         6      +3.320                  +0.238
21023 0 days 03:35:08.336000     11         5      +3.320                  +0.238

[21024 rows x 5 columns]
'@
$textC2 = @'
[datetime.timedelta(0), datetime.timedelta(days=1), datetime.timedelta(days=1)]
'@

$ws.Range("A2").Value = $textA2
$ws.Range("B2").Value = $textB2
$ws.Range("C2").Value = $textC2

# The old rows 3 and 4 are now empty (their content was consolidated into
# row 2 above); clear them so the used range shrinks to A1:C2.
$ws.Range("A3:B4").Clear()
